$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") - copy H1's style (bold, border,
# centered) so the new header cells match the existing ones.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-10 for the new columns I (I0) and J (IF)
$values = @(
    @(9, 9),
    @(9, 9),
    @(6, 7),
    @(7, 8),
    @(8, 8),
    @(8, 9),
    @(9, 9),
    @(6, 6),
    @(5, 5)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
